$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PhanCongNganHan")

# Set the new "Kết Quả Cuối Cùng" (final result) values for the
# QuanLyTrongTai related rows (3, 4, 5) following the 100% integration.
$ws.Range("E4").Value = "100% (9/06/2010)"
$ws.Range("E3").Value = "100% (8/06/2010)"
$ws.Range("D5").Value = "70% (Xong 2/3 Chức năng)"
$ws.Range("E5").Value = "100% (9/06/2010)"

# Update the active selection on the sheet to reflect the latest edit.
$ws.Activate()
$ws.Range("E6").Select()
